$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.161735577932661
$ws.Range("B3").Value = 0.9928586087949914
$ws.Range("B4").Value = 1.48309608044758
$ws.Range("B5").Value = 2.199573983838975
